$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (matches the source file where these
# Price/Volume columns are stored as inline strings, not numbers), then
# restore the default "Normal" style so no stray style index is left on
# the cell (mirrors the original workbook where these cells carry no s=).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) values for rows whose coin identity is unchanged ---
Set-TextValue $ws.Range("D2") '42.087.80'
Set-TextValue $ws.Range("E2") '  -0.06%  '
Set-TextValue $ws.Range("D3") '2.224.42'
Set-TextValue $ws.Range("E3") '  -0.83%  '
Set-TextValue $ws.Range("E4") '  +0.20%  '
Set-TextValue $ws.Range("D5") '243.86'
Set-TextValue $ws.Range("E5") '  -1.29%  '
Set-TextValue $ws.Range("E6") '  +1.57%  '
Set-TextValue $ws.Range("D7") '73.78'
Set-TextValue $ws.Range("E7") '  -0.44%  '
Set-TextValue $ws.Range("E8") '  +0.15%  '
Set-TextValue $ws.Range("E9") '  -0.56%  '
Set-TextValue $ws.Range("D10") '43.38'
Set-TextValue $ws.Range("E10") '  +6.19%  '
Set-TextValue $ws.Range("D11") '0.0963'
Set-TextValue $ws.Range("E11") '  +2.48%  '
Set-TextValue $ws.Range("E12") '  +0.57%  '
Set-TextValue $ws.Range("D13") '0.103'
Set-TextValue $ws.Range("E13") '  +0.48%  '
Set-TextValue $ws.Range("D14") '14.27'
Set-TextValue $ws.Range("E14") '  -1.28%  '
Set-TextValue $ws.Range("D15") '0.844'
Set-TextValue $ws.Range("E15") '  -0.56%  '
Set-TextValue $ws.Range("D16") '2.226.76'
Set-TextValue $ws.Range("E16") '  -0.50%  '
Set-TextValue $ws.Range("D17") '41.976.67'
Set-TextValue $ws.Range("E17") '  -0.05%  '
Set-TextValue $ws.Range("D18") '0.0000110'
Set-TextValue $ws.Range("E18") '  +12.60%  '
Set-TextValue $ws.Range("D19") '6.24'
Set-TextValue $ws.Range("E19") '  +2.08%  '
Set-TextValue $ws.Range("D20") '72.24'
Set-TextValue $ws.Range("E20") '  +0.68%  '
Set-TextValue $ws.Range("D21") '10.22'
Set-TextValue $ws.Range("E21") '  +33.51%  '
Set-TextValue $ws.Range("E22") '  +0.32%  '
Set-TextValue $ws.Range("E23") '  -8.35%  '
Set-TextValue $ws.Range("D24") '11.59'
Set-TextValue $ws.Range("E24") '  +5.20%  '
Set-TextValue $ws.Range("E25") '  -0.01%  '
Set-TextValue $ws.Range("E26") '  +1.12%  '
Set-TextValue $ws.Range("E27") '  -0.50%  '
Set-TextValue $ws.Range("E28") '  +3.00%  '
Set-TextValue $ws.Range("D29") '166.86'
Set-TextValue $ws.Range("E29") '  -1.56%  '
Set-TextValue $ws.Range("D30") '20.65'
Set-TextValue $ws.Range("E30") '  +0.00%  '
Set-TextValue $ws.Range("D31") '5.58'
Set-TextValue $ws.Range("E31") '  +14.67%  '
Set-TextValue $ws.Range("E32") '  -3.38%  '
Set-TextValue $ws.Range("E33") '  +1.11%  '

# --- Rows 34/35: InjectiveProtocol and Kaspa swap positions (rank/index stays, content moves) ---
Set-TextValue $ws.Range("B34") 'InjectiveProtocol'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D34") '29.50'
Set-TextValue $ws.Range("E34") '  -1.36%  '
Set-TextValue $ws.Range("B35") 'Kaspa'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D35") '0.114'
Set-TextValue $ws.Range("E35") '  -3.71%  '

Set-TextValue $ws.Range("E36") '  -4.66%  '
Set-TextValue $ws.Range("E37") '  +0.99%  '
Set-TextValue $ws.Range("D38") '13.09'
Set-TextValue $ws.Range("E38") '  -2.11%  '
Set-TextValue $ws.Range("E39") '  -1.48%  '

# --- Rows 40/41: THORChain and MultiversX swap positions (rank/index stays, content moves) ---
Set-TextValue $ws.Range("B40") 'THORChain'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D40") '5.67'
Set-TextValue $ws.Range("E40") '  -1.72%  '
Set-TextValue $ws.Range("B41") 'MultiversX'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range("D41") '64.46'
Set-TextValue $ws.Range("E41") '  +6.17%  '

Set-TextValue $ws.Range("E42") '  +0.43%  '
Set-TextValue $ws.Range("D43") '104.85'
Set-TextValue $ws.Range("E43") '  -3.35%  '
Set-TextValue $ws.Range("E44") '  +0.72%  '
Set-TextValue $ws.Range("E45") '  +1.04%  '
Set-TextValue $ws.Range("E46") '  +6.48%  '
Set-TextValue $ws.Range("E47") '  +0.49%  '
Set-TextValue $ws.Range("E48") '  +0.86%  '
Set-TextValue $ws.Range("E49") '  +0.77%  '
Set-TextValue $ws.Range("D50") '4.04'
Set-TextValue $ws.Range("E50") '  -0.72%  '
Set-TextValue $ws.Range("D51") '2.432.32'
Set-TextValue $ws.Range("E51") '  -0.87%  '
